$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 52, where E52/F52 hold the
# "Accuracy over PyType" summary line. We need to:
#   1. Fill in C52/D52 with a new "Scalpel Accuracy:" / 1125 summary line
#      (leaving A52/B52/E52/F52 of row 52 blank).
#   2. Push the existing "Accuracy over PyType" / 150 line down into a new
#      row 53 (in E53/F53), clearing it from row 52, while keeping the
#      whole row formatted the same (white fill) as the rest of the table.

$ws.Range("C52").Value = "Scalpel Accuracy:"
$ws.Range("D52").Value = 1125

$ws.Range("A53").Value = ""
$ws.Range("B53").Value = ""
$ws.Range("C53").Value = ""
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = "Accuracy over PyType"
$ws.Range("F53").Value = 150

# Match the existing row formatting (white fill, same as row 52) so the
# new row's cells carry the same style as the rest of the table.
$ws.Range("A53:F53").Interior.Color = 16777215

$ws.Range("E52").Value = $null
$ws.Range("F52").Value = $null
